$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 currently holds "R40"; change its value to the text "1" (a new shared string),
# not the number 1. A leading apostrophe forces text, matching Excel's "quote prefix".
$ws.Range("B11").Value = "'1"
